# Sprint3_Backlog.xlsx — "Started admin assigning lawyer user story"
#
# 1) Story Points edit on sheet row 5 (Unique ID 3.2): 13 -> 8
# 2) Story Points edit on sheet row 17 (Unique ID 7):  1 -> 3
# 3) New backlog item appended as row 26 (Unique ID 16):
#      "As a lawyer I can fill a form" / 3 points / Souidan / Abougabal / Naka
# 4) View state: scroll down a bit and leave the cursor on C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- story point edits -----------------------------------------------------
$ws.Range("C5").Value = 8
$ws.Range("C17").Value = 3

# --- new row -----------------------------------------------------------
$ws.Range("A26").Value = 16
$ws.Range("B26").Value = "As a lawyer I can fill a form"
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = "Souidan"
$ws.Range("E26").Value = "Abougabal"
$ws.Range("F26").Value = "Naka"

# --- view state: scroll the window and land the selection on C18 -----------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$null = $ws.Range("C18").Select()
